# chore: update Sheets via scheduled runner
# Updates recomputed price/profit figures (columns H-N) on a handful of
# leve rows across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets to reflect the
# latest market-board pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 579
$ws.Range("I19").Value = 504.5625
$ws.Range("J19").Value = 678.25
$ws.Range("K19").Value = 504.5625
$ws.Range("L19").Value = 678.25
$ws.Range("M19").Value = -329.5625
$ws.Range("N19").Value = -1028.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 71.84999999999999
$ws.Range("I33").Value = 70.26316
$ws.Range("J33").Value = 102
$ws.Range("K33").Value = 70.26316
$ws.Range("L33").Value = 102
$ws.Range("M33").Value = 158.73684
$ws.Range("N33").Value = -560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 781.35297
$ws.Range("I110").Value = 755.8333
$ws.Range("J110").Value = 842.6
$ws.Range("K110").Value = 755.8333
$ws.Range("L110").Value = 842.6
$ws.Range("M110").Value = 1289.1667
$ws.Range("N110").Value = -4932.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 11112651
$ws.Range("I122").Value = 13334880
$ws.Range("J122").Value = 1504.6666
$ws.Range("K122").Value = 40004640
$ws.Range("L122").Value = 4513.9998
$ws.Range("M122").Value = -40002190
$ws.Range("N122").Value = -9413.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2452533
$ws.Range("I132").Value = 3788707.8
$ws.Range("K132").Value = 11366123.4
$ws.Range("M132").Value = -11363593.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8750.069
$ws.Range("I31").Value = 1223.2069
$ws.Range("J31").Value = 24341.428
$ws.Range("K31").Value = 1223.2069
$ws.Range("L31").Value = 24341.428
$ws.Range("M31").Value = -928.2068999999999
$ws.Range("N31").Value = -24931.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8750.069
$ws.Range("I34").Value = 1223.2069
$ws.Range("J34").Value = 24341.428
$ws.Range("K34").Value = 1223.2069
$ws.Range("L34").Value = 24341.428
$ws.Range("M34").Value = -1021.2069
$ws.Range("N34").Value = -24745.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 496.375
$ws.Range("I5").Value = 496.375
$ws.Range("K5").Value = 1489.125
$ws.Range("M5").Value = -1377.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 7407643.5
$ws.Range("I10").Value = 11111165
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 33333495
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = -33333356
$ws.Range("N10").Value = -2078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 1831.6666
$ws.Range("J104").Value = 1600
$ws.Range("L104").Value = 4800
$ws.Range("N104").Value = -10042

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 490.125
$ws.Range("I122").Value = 388.7857
$ws.Range("J122").Value = 1199.5
$ws.Range("K122").Value = 3499.0713
$ws.Range("L122").Value = 10795.5
$ws.Range("M122").Value = -1049.0713
$ws.Range("N122").Value = -15695.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 496.375
$ws.Range("I135").Value = 496.375
$ws.Range("K135").Value = 4467.375
$ws.Range("M135").Value = -1932.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3160
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3266.6667
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3266.6667
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5262.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3160
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3266.6667
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 16333.3335
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -26317.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 83335510
$ws.Range("I97").Value = 2766.8572
$ws.Range("J97").Value = 200001340
$ws.Range("K97").Value = 2766.8572
$ws.Range("L97").Value = 200001340
$ws.Range("M97").Value = -2270.8572
$ws.Range("N97").Value = -200002332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2600.125
$ws.Range("I122").Value = 1900.5
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 5701.5
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -3251.5
$ws.Range("N122").Value = -13399.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1690.7742
$ws.Range("I126").Value = 1127
$ws.Range("J126").Value = 1959.238
$ws.Range("K126").Value = 3381
$ws.Range("L126").Value = 5877.714
$ws.Range("M126").Value = -911
$ws.Range("N126").Value = -10817.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 71430660
$ws.Range("I132").Value = 95239770
$ws.Range("J132").Value = 3328.2856
$ws.Range("K132").Value = 285719310
$ws.Range("L132").Value = 9984.856800000001
$ws.Range("M132").Value = -285716780
$ws.Range("N132").Value = -15044.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4485.2666
$ws.Range("I7").Value = 3091.3572
$ws.Range("J7").Value = 24000
$ws.Range("K7").Value = 3091.3572
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = -2979.3572
$ws.Range("N7").Value = -24224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 93112.17999999999
$ws.Range("I93").Value = 1666.6666
$ws.Range("J93").Value = 127404.25
$ws.Range("K93").Value = 1666.6666
$ws.Range("L93").Value = 127404.25
$ws.Range("M93").Value = -418.6666
$ws.Range("N93").Value = -129900.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 83337096
$ws.Range("I122").Value = 200002200
$ws.Range("J122").Value = 4871.4287
$ws.Range("K122").Value = 600006600
$ws.Range("L122").Value = 14614.2861
$ws.Range("M122").Value = -600004150
$ws.Range("N122").Value = -19514.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4485.2666
$ws.Range("I126").Value = 3091.3572
$ws.Range("J126").Value = 24000
$ws.Range("K126").Value = 9274.071599999999
$ws.Range("L126").Value = 72000
$ws.Range("M126").Value = -6804.071599999999
$ws.Range("N126").Value = -76940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6367.3335
$ws.Range("I62").Value = 4333.3335
$ws.Range("J62").Value = 7384.3335
$ws.Range("K62").Value = 4333.3335
$ws.Range("L62").Value = 7384.3335
$ws.Range("M62").Value = -3709.3335
$ws.Range("N62").Value = -8632.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6367.3335
$ws.Range("I65").Value = 4333.3335
$ws.Range("J65").Value = 7384.3335
$ws.Range("K65").Value = 21666.6675
$ws.Range("L65").Value = 36921.6675
$ws.Range("M65").Value = -18546.6675
$ws.Range("N65").Value = -43161.6675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 23900
$ws.Range("J75").Value = 23900
$ws.Range("L75").Value = 23900
$ws.Range("N75").Value = -25772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 23900
$ws.Range("J78").Value = 23900
$ws.Range("L78").Value = 71700
$ws.Range("N78").Value = -81060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 83333950
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 125000424
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 375001272
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -375005112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1401.8
$ws.Range("I122").Value = 1952
$ws.Range("J122").Value = 1035
$ws.Range("K122").Value = 5856
$ws.Range("L122").Value = 3105
$ws.Range("M122").Value = -3406
$ws.Range("N122").Value = -8005
